# Update "想去人数" (attendance count) figures in both the "展览" and
# "全部类型" worksheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 82
$wsExpo.Range("F5").Value = 2701
$wsExpo.Range("F8").Value = 7
$wsExpo.Range("F9").Value = 1427
$wsExpo.Range("F13").Value = 1215
$wsExpo.Range("F17").Value = 41
$wsExpo.Range("F18").Value = 38
$wsExpo.Range("F20").Value = 74
$wsExpo.Range("F22").Value = 2625
$wsExpo.Range("F23").Value = 42
$wsExpo.Range("F24").Value = 302

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 82
$wsAll.Range("F5").Value = 2701
$wsAll.Range("F8").Value = 7
$wsAll.Range("F9").Value = 1427
$wsAll.Range("F13").Value = 1215
$wsAll.Range("F17").Value = 42
$wsAll.Range("F18").Value = 38
$wsAll.Range("F20").Value = 74
$wsAll.Range("F22").Value = 2625
$wsAll.Range("F23").Value = 42
$wsAll.Range("F24").Value = 302
